$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly price row above the current row 233; this shifts
# every subsequent row (233..282) down by one, which matches the diff
# (each existing row's contents move to row+1, and the former last row 282
# now also appears, shifted, as row 283).
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with the new record.
$ws.Range("A233").Value = 4
$ws.Range("B233").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C233").Value = "Los Lagos"
$ws.Range("D233").Value = 44782
$ws.Range("E233").Value = 10
$ws.Range("F233").Value = "Fruta"
$ws.Range("G233").Value = 100108
$ws.Range("H233").Value = "Tropicales y subtropicales"
$ws.Range("I233").Value = 100108005
$ws.Range("J233").Value = "Piña"
$ws.Range("K233").Value = "Caramelo"
$ws.Range("L233").Value = "Primera"
$ws.Range("M233").Value = 100
$ws.Range("N233").Value = 23000
$ws.Range("O233").Value = 23000
$ws.Range("P233").Value = 23000
$ws.Range("Q233").Value = "$/caja 12 unidades"
$ws.Range("R233").Value = "Ecuador"
$ws.Range("S233").Value = 1917
$ws.Range("T233").Value = 12
